$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels in row 1 to friendlier display text
$ws.Range("A1").Value = "Admission Number"
$ws.Range("B1").Value = "Symbol Number/Roll Number"
$ws.Range("C1").Value = "Student Name"
$ws.Range("D1").Value = "Date Of Birth(BS)"
$ws.Range("E1").Value = "Religion"
$ws.Range("F1").Value = "Mobile No"
$ws.Range("G1").Value = "Email"
$ws.Range("H1").Value = "Admission Date"
$ws.Range("I1").Value = "Blood Group (O+, A+, B+, AB+, O-, A-, B-, AB-)"
$ws.Range("J1").Value = "Gender(Male/Female)"

# Update view: scroll/selection moved one column to the right
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("J1").Select()
